$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (objective), C (gap), D (solve time), E (Status)
$data = @(
    @{Row=2;  B=-274.9655845122987;   C=0.09045029071781731;  D=3255.784003795; E="OPTIMAL"}
    @{Row=3;  B=-273.98162095885345;  C=3.372166939227947;    D=3603.986717974; E="TIME_LIMIT"}
    @{Row=4;  B=-274.0896045963642;   C=3.2927143918754673;   D=3685.793163656; E="TIME_LIMIT"}
    @{Row=5;  B=-276.8685496566217;   C=4.176958512126726;    D=3624.980113498; E="TIME_LIMIT"}
    @{Row=6;  B=-272.16539220117124;  C=3.6729305785275232;   D=3600.844455122; E="TIME_LIMIT"}
    @{Row=7;  B=-268.97221193176233;  C=0.014916014855304804; D=2155.165336923; E="OPTIMAL"}
    @{Row=8;  B=-265.4281513734784;   C=0.2929961414468023;   D=3769.920243309; E="TIME_LIMIT"}
    @{Row=9;  B=-274.2017067884772;   C=5.69475254826527;     D=3615.965366232; E="TIME_LIMIT"}
    @{Row=10; B=-271.53604103234676;  C=0.22802400293768618;  D=3851.960560323; E="TIME_LIMIT"}
    @{Row=11; B=-268.7867634966758;   C=0.09885584246248551;  D=2597.431866718; E="OPTIMAL"}
)

foreach ($row in $data) {
    $r = $row.Row
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
}
